# Auto-generated Excel COM-interop script to update cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.266.93"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.645.68"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'595.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "'158.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.50%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.538"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").Value = "2.645.39"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").Value = "'0.138"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.07%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "'27.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "3.131.09"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "'0.0000185"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.29%  "
$ws.Range("D17").Value = "68.190.37"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "2.700.98"
$ws.Range("E18").Value = "  +2.39%  "
$ws.Range("D19").Value = "'11.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").Value = "'362.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "'7.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").Value = "'4.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").Value = "'4.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("D24").Value = "'2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("D25").Value = "'75.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "'9.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.22%  "
$ws.Range("D28").Value = "2.785.52"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "'0.0000101"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.54%  "
$ws.Range("D31").Value = "'575.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("D32").Value = "'8.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "'1.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("E35").Value = "  +4.14%  "
$ws.Range("D36").Value = "'0.129"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "'160.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").Value = "'19.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("D40").Value = "'0.371"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").Value = "'1.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "'5.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").Value = "'2.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₆0315"
$ws.Range("E45").Value = "  -7.11%  "
$ws.Range("D46").Value = "'158.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("D47").Value = "'3.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.24%  "
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'0.590"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.46%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'21.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("E51").Value = "  -0.74%  "
